$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Insert a new row at position 13 (pushes the old rows 13-21 down to
# 14-22). This new row carries the professor name that used to (wrongly)
# sit in row 10's B/C cells.
# ---------------------------------------------------------------------
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Clear()
$ws.Range("B10:C10").Copy($ws.Range("B13"))
$ws.Range("B13").Value = "198273 - Domingos Savio Giordani"
$ws.Range("C13").Value = "198273 - Domingos Savio Giordani"

# ---------------------------------------------------------------------
# Row 10 "Objetivos:" now gets its real (pt-br) objectives text instead
# of the professor name.
# ---------------------------------------------------------------------
$ws.Range("B10").Value = "Levar os alunos a vivenciarem de forma mais aprofundada problemas reais da indústria para, em equipes, apresentarem as possíveis soluções, de forma que, com isso, desenvolvam habilidades transversais fundamentais para sua vida profissional, tais como trabalho em equipe, gerenciamento de projetos, pro atividade, ao mesmo tempo em que consolidam o conhecimento adquirido durante o curso."
$ws.Range("C10").Value = "Levar os alunos a vivenciarem de forma mais aprofundada problemas reais da indústria para, em equipes, apresentarem as possíveis soluções, de forma que, com isso, desenvolvam habilidades transversais fundamentais para sua vida profissional, tais como trabalho em equipe, gerenciamento de projetos, pro atividade, ao mesmo tempo em que consolidam o conhecimento adquirido durante o curso."

# ---------------------------------------------------------------------
# Row 14 "Programa resumido:" gets the real pt-br short syllabus text
# instead of the placeholder "Semestral".
# ---------------------------------------------------------------------
$ws.Range("B14").Value = "Formação e trabalho em equipes, Comunicação, Inovação Sistemática, Legislação, Gerenciamento de Projetos. Identificação de Problemas, Formulação do Projeto, Especificação de Problemas, Análise do Conhecimento disponível, Avaliação e Tomada de Decisão, Cronograma, Elaboração de relatórios, Apresentação de Projetos"
$ws.Range("C14").Value = "Formação e trabalho em equipes, Comunicação, Inovação Sistemática, Legislação, Gerenciamento de Projetos. Identificação de Problemas, Formulação do Projeto, Especificação de Problemas, Análise do Conhecimento disponível, Avaliação e Tomada de Decisão, Cronograma, Elaboração de relatórios, Apresentação de Projetos"

# ---------------------------------------------------------------------
# Row 16 "Programa:" gets the full pt-br syllabus text instead of the
# stray "01/01/2020" date value.
# ---------------------------------------------------------------------
$ws.Range("B16").Value = "Formação e trabalho em equipes e Comunicação – o desenvolvimento das habilidades essenciais para o trabalho em equipes; Inovação Sistemática – desenvolvimento de soluções inovadoras, sistematização e características; Legislação - noções da legislação aplicada à ação empresarial; Gerenciamento de Projetos e Cronograma – Metodologias e esquematizações necessárias com os elementos gerenciais; Identificação de Problemas – sistematização de ações para a localização de causas; Formulação do Projeto – apresentação dos aspectos gerenciais necessários ao desenvolvimento do projeto, Plano de gestão, Estrutura Analítica do Projeto (EAP) etc; Especificação de Problemas – sistematização dos problemas dentro das áreas de conhecimento; Análise do Conhecimento disponível, Avaliação e Tomada de Decisão; Elaboração de relatórios – formatação dentro das normas ABNT; Apresentação de Projetos."
$ws.Range("C16").Value = "Formação e trabalho em equipes e Comunicação – o desenvolvimento das habilidades essenciais para o trabalho em equipes; Inovação Sistemática – desenvolvimento de soluções inovadoras, sistematização e características; Legislação - noções da legislação aplicada à ação empresarial; Gerenciamento de Projetos e Cronograma – Metodologias e esquematizações necessárias com os elementos gerenciais; Identificação de Problemas – sistematização de ações para a localização de causas; Formulação do Projeto – apresentação dos aspectos gerenciais necessários ao desenvolvimento do projeto, Plano de gestão, Estrutura Analítica do Projeto (EAP) etc; Especificação de Problemas – sistematização dos problemas dentro das áreas de conhecimento; Análise do Conhecimento disponível, Avaliação e Tomada de Decisão; Elaboração de relatórios – formatação dentro das normas ABNT; Apresentação de Projetos."

# ---------------------------------------------------------------------
# Row 19 "Método:" gets the real evaluation-method text instead of the
# professor name that had ended up there.
# ---------------------------------------------------------------------
$ws.Range("B19").Value = "Apresentações intermediárias e finais."
$ws.Range("C19").Value = "Apresentações intermediárias e finais."

# ---------------------------------------------------------------------
# Row 20 "Critério:" gets the real grading-criteria text.
# ---------------------------------------------------------------------
$ws.Range("B20").Value = "Serão feitas duas avaliações por uma banca de professores que assistirão às apresentações, as notas serão as médias das notas dadas pelos professores."
$ws.Range("C20").Value = "Serão feitas duas avaliações por uma banca de professores que assistirão às apresentações, as notas serão as médias das notas dadas pelos professores."

# ---------------------------------------------------------------------
# Row 21 "Norma de recuperação:" gets the real recovery-norm text.
# ---------------------------------------------------------------------
$ws.Range("B21").Value = "Reapresentação do último seminário, cuja nota constituirá a nota final da disciplina."
$ws.Range("C21").Value = "Reapresentação do último seminário, cuja nota constituirá a nota final da disciplina."

# ---------------------------------------------------------------------
# Row 22 "Bibliografia:" gets the real bibliography text instead of the
# recovery-norm text that had been misplaced there.
# ---------------------------------------------------------------------
$ws.Range("B22").Value = "Gestão de Negócios: Visões e dimensões empresariais da o Organização. Autores: Cruz Jr, J.B., Rocha, J.A.O. e Tachizawa, T.Editora: ATLASGestão Empresarial - de Taylor aos nossos diasAutores: Pereira, M. I. , Autor: Ferreira, A. A. e Reis, A.C. F Editora: THOMSON PIONEIRABaron e Shane: Empreendedorismo: uma visão do processo (EVP), Ed. Thomson, 2006Textos fornecidos pelo professor da disciplinaArtigos extraídos de revistas especializadas na área de gestão e produção."
$ws.Range("C22").Value = "Gestão de Negócios: Visões e dimensões empresariais da o Organização. Autores: Cruz Jr, J.B., Rocha, J.A.O. e Tachizawa, T.Editora: ATLASGestão Empresarial - de Taylor aos nossos diasAutores: Pereira, M. I. , Autor: Ferreira, A. A. e Reis, A.C. F Editora: THOMSON PIONEIRABaron e Shane: Empreendedorismo: uma visão do processo (EVP), Ed. Thomson, 2006Textos fornecidos pelo professor da disciplinaArtigos extraídos de revistas especializadas na área de gestão e produção."
